# Update the "Förändrad" (Changed) date column (C) for rows 2-20
# from 2023-10-09 (45208) to 2023-10-13 (45212) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
